# feat: add 2022-Q4 data
#
# The workbook gains a new "2022-Q4" detail sheet (inserted right after the
# "总计" summary sheet, before the existing "2022-Q3" sheet), and the "总计"
# summary sheet gets a new top row for 2022-Q4 with the existing 2022-Q3 /
# 2022-Q2 rows shifting down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating the "2022-Q3" sheet
#    (keeps header styling/borders identical) and placing it right before
#    that source sheet, then trim it down to the rows Q4 actually needs.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
# Copy(Before:=$q3) drops the duplicate immediately in front of $q3, i.e. at
# sheet index 2 (总计 is index 1) - $q3.Index itself isn't refreshed after
# the copy, so address the new sheet positionally instead.
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The Q3 sheet has 5 fund rows (rows 2-6); Q4 only has 3 funds, so drop the
# two extra rows.
$q4.Rows.Item(5).Delete()
$q4.Rows.Item(5).Delete()

# Columns B:G hold text data (fund code/name/size/position/ratio/value) even
# when the text looks numeric (e.g. "000974", "2.33") - force text format
# before assigning so leading zeros / exact text are preserved.
$q4.Range("B1:H1").NumberFormat = "@"
$q4.Range("B2:G4").NumberFormat = "@"

# Header row (already copied from Q3, but keep explicit for clarity/safety)
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Row 2 - 安信消费医药主题股票
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "000974"
$q4.Cells.Item(2,3).Value = "安信消费医药主题股票"
$q4.Cells.Item(2,4).Value = "2.33"
$q4.Cells.Item(2,5).Value = "92.45"
$q4.Cells.Item(2,6).Value = "5.18"
$q4.Cells.Item(2,7).Value = "0.1207"
$q4.Cells.Item(2,8).Value = 3

# Row 3 - 合煦智远消费主题股票C
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "007288"
$q4.Cells.Item(3,3).Value = "合煦智远消费主题股票C"
$q4.Cells.Item(3,4).Value = "0.11"
$q4.Cells.Item(3,5).Value = "83.65"
$q4.Cells.Item(3,6).Value = "4.19"
$q4.Cells.Item(3,7).Value = "0.0046"
$q4.Cells.Item(3,8).Value = 5

# Row 4 - 合煦智远消费主题股票A
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "007287"
$q4.Cells.Item(4,3).Value = "合煦智远消费主题股票A"
$q4.Cells.Item(4,4).Value = "0.03"
$q4.Cells.Item(4,5).Value = "83.65"
$q4.Cells.Item(4,6).Value = "4.19"
$q4.Cells.Item(4,7).Value = "0.0013"
$q4.Cells.Item(4,8).Value = 5

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert the 2022-Q4 figures as the new
#    row 2, pushing 2022-Q3 to row 3 and 2022-Q2 to row 4.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Shift the old rows down first (read old values before they are overwritten).
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q2"
$total.Cells.Item(4,3).Value = 2
$total.Cells.Item(4,4).Value = 0.04

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q3"
$total.Cells.Item(3,3).Value = 5
$total.Cells.Item(3,4).Value = 0.16

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 3
$total.Cells.Item(2,4).Value = 0.13

# New A4 index cell needs the same bold/border/center styling as A2/A3 - copy
# it from A3 (same-sheet format copy keeps the existing style index intact).
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)
$total.Cells.Item(4,1).Value = 2

# ---------------------------------------------------------------------
# 3. Keep the originally-selected tab ("2022-Q2") active - duplicating the
#    Q3 sheet leaves the new copy active, so restore the prior selection.
# ---------------------------------------------------------------------
$wb.Worksheets.Item(4).Activate()
